# Insert a new daily price record for Cilantro (Femacal de La Calera) at
# row 584. Excel shifts every existing row from 584..643 down to 585..644
# (keeping their original cell values/formatting), and the freshly
# inserted row 584 is populated with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(584).Insert()

$ws.Cells.Item(584, 1).Value  = 3
$ws.Cells.Item(584, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(584, 3).Value  = "Coquimbo"
$ws.Cells.Item(584, 4).Value  = 45194
$ws.Cells.Item(584, 5).Value  = 5
$ws.Cells.Item(584, 6).Value  = 100112040
$ws.Cells.Item(584, 7).Value  = "Cilantro"
$ws.Cells.Item(584, 8).Value  = "Sin especificar"
$ws.Cells.Item(584, 9).Value  = "Primera"
$ws.Cells.Item(584, 10).Value = 120
$ws.Cells.Item(584, 11).Value = 4000
$ws.Cells.Item(584, 12).Value = 4000
$ws.Cells.Item(584, 13).Value = 4000
$ws.Cells.Item(584, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(584, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(584, 16).Value = 1333
$ws.Cells.Item(584, 17).Value = 3
$ws.Cells.Item(584, 18).Value = "Hortaliza"
